$d = $word.ActiveDocument

# 1. Title: "hamza arslan" -> "Hamza Arslan"
$r = $d.Paragraphs(1).Range
$r.Find.Execute("hamza arslan", $true, $false, $false, $false, $false,
                $true, 1, $false, "Hamza Arslan", 2)

# 2. Personal-info paragraph: 6 bold runs, each "Label: value", separated by <w:br/>.
#    Replace every value in turn.
$r = $d.Paragraphs(2).Range
$r.Find.Execute("Cinsiyet: erkek", $true, $false, $false, $false, $false,
                $true, 1, $false, "Cinsiyet: Erkek", 2)

$r = $d.Paragraphs(2).Range
$r.Find.Execute("Doğum tarihi: 11.11.1111", $true, $false, $false, $false, $false,
                $true, 1, $false, "Doğum tarihi: 05.11.2000", 2)

$r = $d.Paragraphs(2).Range
$r.Find.Execute("Medeni durumu: bekar", $true, $false, $false, $false, $false,
                $true, 1, $false, "Medeni durumu: Bekar", 2)

$r = $d.Paragraphs(2).Range
$r.Find.Execute("Ülkesi: tr", $true, $false, $false, $false, $false,
                $true, 1, $false, "Ülkesi: TR", 2)

$r = $d.Paragraphs(2).Range
$r.Find.Execute("Askerlik durumu: done", $true, $false, $false, $false, $false,
                $true, 1, $false, "Askerlik durumu: Yapıldı", 2)

$r = $d.Paragraphs(2).Range
$r.Find.Execute("Ehliyet türü: b", $true, $false, $false, $false, $false,
                $true, 1, $false, "Ehliyet türü: B", 2)

# The Find/Replace calls above coalesce neighbouring runs that end up with
# identical formatting (same <w:rPr>) after being touched back-to-back, which
# collapses the original 6-run paragraph into fewer runs. The source document
# keeps each "Label: value" segment (and its leading <w:br/>) in its own run,
# so re-split the paragraph at each original run boundary by toggling Bold
# off/on across the break + following text - this forces the engine to emit
# a fresh run starting at that break, restoring the original run layout.
$labels = @("Doğum tarihi: 05.11.2000", "Medeni durumu: Bekar", "Ülkesi: TR", "Askerlik durumu: Yapıldı", "Ehliyet türü: B")
foreach ($t in $labels) {
    $f = $d.Paragraphs(2).Range
    $f.Find.Execute($t, $true, $false, $false, $false, $false,
                    $true, 1, $false, "", 0)
    $boundary = $d.Range($f.Start - 1, $f.End)
    $boundary.Bold = 0
    $boundary.Bold = 1
}
